$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: H12 and I12 get value 5; clear J12 (formula result recalculated elsewhere)
$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 5
$ws.Range("J12").ClearContents()

# Row 18: C18, D18, E18 change from 2 to 5
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 5

# Rows 20, 21, 22: G column gets value 5
$ws.Range("G20").Value = 5
$ws.Range("G21").Value = 5
$ws.Range("G22").Value = 5

# Update selection to F18
$ws.Range("F18").Select()
